# Append a duplicate copy of the existing project_index data rows
# (A2:A139) to the end of the sheet (A140:A277), matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcStart = 2
$srcEnd = 139
$destStart = 140

$count = $srcEnd - $srcStart + 1

for ($i = 0; $i -lt $count; $i++) {
    $srcRowNum = $srcStart + $i
    $destRowNum = $destStart + $i
    $value = $ws.Cells.Item($srcRowNum, 1).Value2
    $ws.Cells.Item($destRowNum, 1).Value = $value
}
